$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: (row, C, D, E, F, G, K, L, M, N, O)
$data = @(
  ,(2, 13.02524186462899, 5.915200010421848, 13.31724460201748, 29.56490511193664, 3.663773078584197, 14.65813805415564, 9.141208926191613, 17.66702338506492, 18.42617641486155, 26.40903883920713)
  ,(3, 12.99237485270437, 5.899631175168965, 13.33486589885283, 29.57313151650409, 3.666122533740872, 14.28602164749289, 9.163151964872025, 17.52263150501424, 18.48109611978033, 26.4687764412782)
  ,(4, 12.97516797331814, 5.88998372741167, 13.34812599058976, 29.58656884987664, 3.667641749796249, 14.05469630478319, 9.177574510438511, 17.43636932123779, 18.51667331269724, 26.5119151841128)
  ,(5, 12.96890880998883, 5.88603107484191, 13.35414286991662, 29.59415043907299, 3.66828017758675, 13.95985243191095, 9.183690895368722, 17.40185059057594, 18.53163917370285, 26.53111397817484)
  ,(6, 12.96791507913725, 5.885373483848872, 13.35517899421632, 29.59553642367025, 3.66838735766643, 13.94407300918116, 9.184720966300771, 17.39615796287288, 18.53415253385818, 26.5343995799776)
  ,(7, 12.97508050597476, 5.889930505355516, 13.34820465381185, 29.58666257674967, 3.667650281486788, 14.0534193553328, 9.177656029667352, 17.43590118242889, 18.51687325171448, 26.51216755568474)
  ,(8, 13.01329540126471, 5.909850715689705, 13.32281368316003, 29.56599958484291, 3.664567301780195, 14.53050349503874, 9.148578034818193, 17.61675923770903, 18.44472811461923, 26.42829283096368)
  ,(9, 13.11157375839522, 5.948175496848202, 13.29240049633209, 29.5920971375758, 3.659126831806954, 15.43778517521276, 9.099074335699649, 17.98899112012525, 18.31793098487695, 26.31527295611915)
  ,(10, 13.19761753572336, 5.975830706916745, 13.28188160515943, 29.65189972766114, 3.655494676827001, 16.08003737259162, 9.067266452593687, 18.27120600679472, 18.23365338422957, 26.26386272190527)
  ,(11, 13.23966793025767, 5.988291752023856, 13.27966356148642, 29.68790031525246, 3.653920708394852, 16.36560062075381, 9.05378258449185, 18.40105803866345, 18.19722731376843, 26.24738729662308)
  ,(12, 13.2560001945702, 5.992992273700875, 13.27919241306568, 29.70279274706497, 3.653335884504668, 16.47268751375622, 9.048818015663523, 18.45040443496225, 18.18370765375074, 26.2421451765791)
  ,(13, 13.25246473175203, 5.991980758368634, 13.27927748951679, 29.69952946323395, 3.653461339422705, 16.44967263182592, 9.049880936189355, 18.43976964961407, 18.18660717731726, 26.24322978860883)
  ,(14, 13.24100346177292, 5.988678837545085, 13.27961741169183, 29.68910031110179, 3.653872370372055, 16.37443230727645, 9.053371312457784, 18.40511451779763, 18.19610955470509, 26.2469360327302)
  ,(15, 13.23403603174607, 5.986653914192723, 13.2798736356768, 29.68287603902145, 3.65412559604455, 16.32820589411664, 9.055527685776839, 18.38390883432878, 18.20196570591105, 26.24933610028892)
  ,(16, 13.19492717012234, 5.975013856798571, 13.28207818821685, 29.6497236059777, 3.655599110329119, 16.06123353085122, 9.068167467276663, 18.26274628646943, 18.23607231729925, 26.26507876557988)
  ,(17, 13.17167352754829, 5.967841733738574, 13.28408789131205, 29.6316354339342, 3.656523081019205, 15.89568980055196, 9.076173836221074, 18.188767256133, 18.2574847573805, 26.27650878955365)
  ,(18, 13.15857303873148, 5.963705336906512, 13.28548547527473, 29.62205992798876, 3.657061899542534, 15.7998576171609, 9.0808716817955, 18.14635717864181, 18.26998066563987, 26.28373342391339)
  ,(19, 13.15418484601296, 5.962302934922064, 13.2860001859979, 29.61896022230148, 3.657245602681241, 15.76730792608519, 9.082478237115609, 18.13202315628154, 18.27424251441708, 26.28629116790863)
  ,(20, 13.17412058625401, 5.968606383143894, 13.28384894657745, 29.63347525324557, 3.656423959896759, 15.91337670279541, 9.075311942671796, 18.19662815780402, 18.25518673960934, 26.27522470849751)
  ,(21, 13.24435890226001, 5.989649194041158, 13.27950756358004, 29.6921294636895, 3.653751337012432, 16.39656146233527, 9.052342266700842, 18.41528913821798, 18.19331104391268, 26.24582034639141)
  ,(22, 13.29264098754763, 6.003295291664204, 13.27881943809935, 29.73780295542715, 3.652069901095459, 16.70619146006565, 9.038154787060995, 18.55919512038723, 18.15446898968679, 26.23241378802473)
  ,(23, 13.26665773465508, 5.996022206741554, 13.27899021523271, 29.71275665977014, 3.652961361267328, 16.54153068981238, 9.045651552895052, 18.48231075124683, 18.17505385544157, 26.23903659613148)
  ,(24, 13.17301343445677, 5.968260725693431, 13.28395621910358, 29.63264090497795, 3.656468748830343, 15.90538249882276, 9.075701309185344, 18.19307386298533, 18.25622509518365, 26.27580320653386)
  ,(25, 13.0825255152478, 5.937892023222382, 13.29855159311502, 29.5778975297392, 3.660534244152554, 15.19613305164122, 9.111663693543891, 17.88661287983776, 18.35066863641369, 26.34031070760599)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 3).Value  = $row[1]
  $ws.Cells.Item($r, 4).Value  = $row[2]
  $ws.Cells.Item($r, 5).Value  = $row[3]
  $ws.Cells.Item($r, 6).Value  = $row[4]
  $ws.Cells.Item($r, 7).Value  = $row[5]
  $ws.Cells.Item($r, 11).Value = $row[6]
  $ws.Cells.Item($r, 12).Value = $row[7]
  $ws.Cells.Item($r, 13).Value = $row[8]
  $ws.Cells.Item($r, 14).Value = $row[9]
  $ws.Cells.Item($r, 15).Value = $row[10]
}
